# Apply targeted edits to the QC sample type CodeSystem workbook.

$wb = $excel.ActiveWorkbook

# --- Update the "Date" value on the Metadata sheet ---
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Range("B8").Value = "2025-08-13T14:50:19+00:00"

# --- Fix the malformed NCIT codes (missing "C" prefix) on the Concepts sheet ---
$conceptsSheet = $wb.Worksheets.Item("Concepts")
$conceptsSheet.Range("B11").Value = "NCIT:C156440"
$conceptsSheet.Range("B12").Value = "NCIT:C156441"
$conceptsSheet.Range("B13").Value = "NCIT:C164032"
